$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 18 (filled C, then B, then A)
$ws.Cells.Item(18, 3).Value = "第68期 混合式紙飛機"
$ws.Cells.Item(18, 2).Value = "9/19"
$ws.Cells.Item(18, 1).Value = "7/25"

# Row 19 (filled B, then C, then A)
$ws.Cells.Item(19, 2).Value = "9/26"
$ws.Cells.Item(19, 3).Value = "第69期 祕寶 開放區域 天吉禍 祕寶效果: 烘培貓貓包有10% (36) 機率麵粉翻倍"
$ws.Cells.Item(19, 1).Value = "8/1"

# Row 20 (filled B, then C, then A)
$ws.Cells.Item(20, 2).Value = "10/3"
$ws.Cells.Item(20, 3).Value = "第70期 第四代寵物"
$ws.Cells.Item(20, 1).Value = "8/8"

# Row 21 (filled B, then C, then A)
$ws.Cells.Item(21, 2).Value = "10/10"
$ws.Cells.Item(21, 3).Value = "第71期 星途 "
$ws.Cells.Item(21, 1).Value = "8/15"

# Match source formatting: column A/B style (numFmtId 49 text, applied to existing A/B cells) for new rows
$ws.Range("A18:B21").NumberFormat = "@"

# Update the view to match the edited file: scrolled down with C21 selected
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 13
$excel.ActiveWindow.ScrollColumn = 1
[void]$ws.Range("C21").Select()
